# OkraFreezing.docx edit
#
# Splits the opening paragraph:
#   "The smooth type varieties ... split as easily. Preparation - Select ... seed cell."
# into three paragraphs, matching the blank-line-separated "method" formatting
# already used by the rest of the document ("Water blanch...", "Cool promptly...", etc.):
#
#   P1: "The smooth type varieties ... split as easily."      (+ trailing literal "\n" run)
#   P2: (empty paragraph, just the literal "\n" run)
#   P3: "Preparation - Select ... seed cell."                 (+ trailing literal "\n" run)

$d = $word.ActiveDocument

# --- Locate the start of "Preparation" (the text that must become its own paragraph) ---
$find1 = $d.Content
$ok1 = $find1.Find.Execute("Preparation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok1) {
    throw "Could not find 'Preparation' in document"
}
$prepStart = $find1.Start

# --- Split #1: break the paragraph right before "Preparation" ---
# (Using InsertParagraphBefore here, rather than editing/inserting into the tail of the
#  paragraph, keeps the pre-existing trailing "\n" run attached to the new "Preparation..."
#  paragraph instead of merging it into a freshly rebuilt run.)
$d.Range($prepStart, $prepStart).InsertParagraphBefore()

# --- Split #2: insert another paragraph break right before "Preparation" again ---
# This produces an empty paragraph between the "...easily." paragraph and the
# "Preparation..." paragraph.
$find2 = $d.Content
$ok2 = $find2.Find.Execute("Preparation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) {
    throw "Could not re-find 'Preparation' after first split"
}
$prepStart2 = $find2.Start
$d.Range($prepStart2, $prepStart2).InsertParagraphBefore()

# --- Fix up paragraph 1: "...split as easily. " -> "...split as easily." + "\n" run ---
# After the two splits above, paragraph 1 ends with a single trailing space
# (the space that used to separate "easily." and "Preparation") right before its
# paragraph mark. Turn that trailing space into the literal text "\n" (matching the
# "\n" marker runs used throughout the rest of the document).
$p1 = $d.Paragraphs(1).Range
$spacePos = $p1.End - 2
$spaceCheck = $d.Range($spacePos, $spacePos + 1).Text
if ($spaceCheck -ne " ") {
    throw "Unexpected character before 'Preparation' split point: [$spaceCheck]"
}
$d.Range($spacePos, $spacePos + 1).Text = "\n"

# --- Fix up paragraph 2: empty paragraph -> literal "\n" text ---
$p2 = $d.Paragraphs(2).Range
$d.Range($p2.Start, $p2.Start).InsertAfter("\n")
